$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# New bitmap-related rows appended after the existing data (row 29),
# leaving row 30 blank just like the rest of the sheet's layout.
$ws.Range("A31").Value2 = "BitmapTestCount"
$ws.Range("B31").Value2 = 7

$ws.Range("A32").Value2 = "AllowDebug"
$ws.Range("B32").Value2 = $True

$ws.Range("A33").Value2 = "ModeSelect"
$ws.Range("B33").Value2 = -3

$ws.Range("A34").Value2 = "RegionCode"
$ws.Range("B34").Value2 = -8

$ws.Range("A35").Value2 = "PowerGood"
$ws.Range("B35").Value2 = $True

$ws.Range("A36").Value2 = "FanRunning"
$ws.Range("B36").Value2 = $False

$ws.Range("A37").Value2 = "ErrorCode"
$ws.Range("B37").Value2 = 42

$ws.Range("A38").Value2 = "HwRevision"
$ws.Range("B38").Value2 = 15

$ws.Range("A39").Value2 = "BitmapCheckVal"
$ws.Range("B39").Value2 = 305419896
$ws.Range("B39").NumberFormat = "#,##0"

# Match the saved selection state recorded in the workbook.
[void]$ws.Range("D36").Select()
